# SM18 microstate list: remove 14 replicate microstates (v1.4.1).
#
# The workbook lists molecular microstates (ID in column B, canonical
# isomeric SMILES in column C, rows 3..67) with one 2-D depiction picture
# anchored per row. This edit removes 14 rows that were identified as
# duplicate/replicate microstates, together with their associated shared
# strings and one picture per removed row (the trailing 14 picture shapes),
# shifting everything else up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Microstate IDs (column B) to remove.
$idsToRemove = @(
    "SM18_micro008",
    "SM18_micro015",
    "SM18_micro019",
    "SM18_micro027",
    "SM18_micro035",
    "SM18_micro039",
    "SM18_micro040",
    "SM18_micro041",
    "SM18_micro043",
    "SM18_micro044",
    "SM18_micro046",
    "SM18_micro066",
    "SM18_micro067",
    "SM18_micro073"
)

# Drop the trailing picture shapes (one per removed row) so the count of
# 2-D depiction images again matches the number of remaining data rows.
$shapeCount = $ws.Shapes.Count
$shapesToDrop = $idsToRemove.Count
for ($i = $shapeCount; $i -ge ($shapeCount - $shapesToDrop + 1); $i--) {
    $ws.Shapes.Item($i).Delete()
}

# Delete the matching data rows, bottom-to-top so row numbers of
# not-yet-processed rows stay valid while we work our way up.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = $lastRow; $r -ge 1; $r--) {
    $microstateId = $ws.Cells.Item($r, 2).Value()
    if ($idsToRemove -contains $microstateId) {
        $ws.Rows("$r`:$r").Delete()
    }
}
